$wb = $excel.ActiveWorkbook

# This script applies a batch of market-data refresh values to the
# "Goblin_Profits" price-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each worksheet tracks, per leve (crafting order), live Market Board prices
# pulled by a scheduled data-refresh job; columns H-N hold those refreshed values.

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 343.52942
$ws.Range("J9").Value = 395.7143
$ws.Range("L9").Value = 395.7143
$ws.Range("N9").Value = -733.7143
$ws.Range("H80").Value = 1919.125
$ws.Range("I80").Value = 739.8570999999999
$ws.Range("J80").Value = 2836.3333
$ws.Range("K80").Value = 2219.5713
$ws.Range("L80").Value = 8508.999899999999
$ws.Range("M80").Value = -1221.5713
$ws.Range("N80").Value = -10504.9999
$ws.Range("H83").Value = 1919.125
$ws.Range("I83").Value = 739.8570999999999
$ws.Range("J83").Value = 2836.3333
$ws.Range("K83").Value = 6658.7139
$ws.Range("L83").Value = 25526.9997
$ws.Range("M83").Value = -1666.7139
$ws.Range("N83").Value = -35510.9997
$ws.Range("H112").Value = 1474.5
$ws.Range("I112").Value = 1051.8
$ws.Range("K112").Value = 3155.4
$ws.Range("M112").Value = -2047.4
$ws.Range("H113").Value = 3603.8262
$ws.Range("J113").Value = 3997.25
$ws.Range("L113").Value = 3997.25
$ws.Range("N113").Value = -10505.25
$ws.Range("H130").Value = 74990
$ws.Range("J130").Value = 74990
$ws.Range("L130").Value = 74990
$ws.Range("N130").Value = -85030
$ws.Range("H132").Value = 1455.826
$ws.Range("I132").Value = 1110.4872
$ws.Range("J132").Value = 3379.8572
$ws.Range("K132").Value = 3331.4616
$ws.Range("L132").Value = 10139.5716
$ws.Range("M132").Value = -801.4616000000001
$ws.Range("N132").Value = -15199.5716
$ws.Range("H137").Value = 2979.9375
$ws.Range("I137").Value = 1985.1428
$ws.Range("J137").Value = 3753.6667
$ws.Range("K137").Value = 5955.428400000001
$ws.Range("L137").Value = 11261.0001
$ws.Range("M137").Value = -3405.428400000001
$ws.Range("N137").Value = -16361.0001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5820.8667
$ws.Range("I2").Value = 9553.5
$ws.Range("J2").Value = 3332.4443
$ws.Range("K2").Value = 9553.5
$ws.Range("L2").Value = 3332.4443
$ws.Range("M2").Value = -9440.5
$ws.Range("N2").Value = -3558.4443
$ws.Range("H61").Value = 7347.033
$ws.Range("I61").Value = 7071.769
$ws.Range("J61").Value = 7557.5293
$ws.Range("K61").Value = 7071.769
$ws.Range("L61").Value = 7557.5293
$ws.Range("M61").Value = -6859.769
$ws.Range("N61").Value = -7981.5293
$ws.Range("H63").Value = 10048.375
$ws.Range("I63").Value = 5750
$ws.Range("J63").Value = 11481.167
$ws.Range("K63").Value = 5750
$ws.Range("L63").Value = 11481.167
$ws.Range("M63").Value = -5064
$ws.Range("N63").Value = -12853.167
$ws.Range("H66").Value = 10048.375
$ws.Range("I66").Value = 5750
$ws.Range("J66").Value = 11481.167
$ws.Range("K66").Value = 28750
$ws.Range("L66").Value = 57405.835
$ws.Range("M66").Value = -25318
$ws.Range("N66").Value = -64269.835
$ws.Range("H102").Value = 8313.477000000001
$ws.Range("I102").Value = 6365.3335
$ws.Range("K102").Value = 6365.3335
$ws.Range("M102").Value = -4743.3335
$ws.Range("H116").Value = 5820.8667
$ws.Range("I116").Value = 9553.5
$ws.Range("J116").Value = 3332.4443
$ws.Range("K116").Value = 9553.5
$ws.Range("L116").Value = 3332.4443
$ws.Range("M116").Value = -7259.5
$ws.Range("N116").Value = -7920.4443
$ws.Range("H136").Value = 7347.033
$ws.Range("I136").Value = 7071.769
$ws.Range("J136").Value = 7557.5293
$ws.Range("K136").Value = 21215.307
$ws.Range("L136").Value = 22672.5879
$ws.Range("M136").Value = -18665.307
$ws.Range("N136").Value = -27772.5879

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5820.8667
$ws.Range("I3").Value = 9553.5
$ws.Range("J3").Value = 3332.4443
$ws.Range("K3").Value = 9553.5
$ws.Range("L3").Value = 3332.4443
$ws.Range("M3").Value = -9439.5
$ws.Range("N3").Value = -3560.4443
$ws.Range("H20").Value = 4531.909
$ws.Range("I20").Value = 5700.3477
$ws.Range("J20").Value = 1844.5
$ws.Range("K20").Value = 5700.3477
$ws.Range("L20").Value = 1844.5
$ws.Range("M20").Value = -5453.3477
$ws.Range("N20").Value = -2338.5
$ws.Range("H26").Value = 14024.286
$ws.Range("I26").Value = 14024.286
$ws.Range("K26").Value = 14024.286
$ws.Range("M26").Value = -13732.286
$ws.Range("H64").Value = 752.1111
$ws.Range("J64").Value = 864
$ws.Range("L64").Value = 864
$ws.Range("N64").Value = -1314
$ws.Range("H67").Value = 752.1111
$ws.Range("J67").Value = 864
$ws.Range("L67").Value = 864
$ws.Range("N67").Value = -2424
$ws.Range("H96").Value = 19751.455
$ws.Range("I96").Value = 13227.1
$ws.Range("J96").Value = 84995
$ws.Range("K96").Value = 13227.1
$ws.Range("L96").Value = 84995
$ws.Range("M96").Value = -10481.1
$ws.Range("N96").Value = -90487
$ws.Range("H105").Value = 37000
$ws.Range("I105").Value = 100000
$ws.Range("K105").Value = 100000
$ws.Range("M105").Value = -98253

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.625
$ws.Range("I7").Value = 377.14285
$ws.Range("J7").Value = 75.77778000000001
$ws.Range("K7").Value = 377.14285
$ws.Range("L7").Value = 75.77778000000001
$ws.Range("M7").Value = -264.14285
$ws.Range("N7").Value = -301.77778
$ws.Range("H22").Value = 1352.5714
$ws.Range("I22").Value = 877.2857
$ws.Range("J22").Value = 1827.8572
$ws.Range("K22").Value = 877.2857
$ws.Range("L22").Value = 1827.8572
$ws.Range("M22").Value = -527.2857
$ws.Range("N22").Value = -2527.8572
$ws.Range("H105").Value = 2619.8667
$ws.Range("I105").Value = 2834.5881
$ws.Range("K105").Value = 2834.5881
$ws.Range("M105").Value = -1087.5881
$ws.Range("H140").Value = 350000
$ws.Range("J140").Value = 350000
$ws.Range("L140").Value = 350000
$ws.Range("N140").Value = -360360
$ws.Range("H141").Value = 235690.94
$ws.Range("J141").Value = 268104.94
$ws.Range("L141").Value = 268104.94
$ws.Range("N141").Value = -278464.94

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 910.25
$ws.Range("I110").Value = 910.25
$ws.Range("K110").Value = 2730.75
$ws.Range("M110").Value = 1359.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 253747
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 253747
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = -253979
$ws.Range("H10").Value = 260000
$ws.Range("J10").Value = 260000
$ws.Range("L10").Value = 260000
$ws.Range("N10").Value = -260338
$ws.Range("H80").Value = 5897.8887
$ws.Range("I80").Value = 4243.1816
$ws.Range("J80").Value = 7035.5
$ws.Range("K80").Value = 4243.1816
$ws.Range("L80").Value = 7035.5
$ws.Range("M80").Value = -3245.1816
$ws.Range("N80").Value = -9031.5
$ws.Range("H83").Value = 5897.8887
$ws.Range("I83").Value = 4243.1816
$ws.Range("J83").Value = 7035.5
$ws.Range("K83").Value = 21215.908
$ws.Range("L83").Value = 35177.5
$ws.Range("M83").Value = -16223.908
$ws.Range("N83").Value = -45161.5
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344
$ws.Range("H122").Value = 4558.552
$ws.Range("I122").Value = 3499.5
$ws.Range("J122").Value = 5862
$ws.Range("K122").Value = 10498.5
$ws.Range("L122").Value = 17586
$ws.Range("M122").Value = -8048.5
$ws.Range("N122").Value = -22486
$ws.Range("H132").Value = 2051.875
$ws.Range("I132").Value = 1816.0555
$ws.Range("J132").Value = 2163.5789
$ws.Range("K132").Value = 5448.166499999999
$ws.Range("L132").Value = 6490.736699999999
$ws.Range("M132").Value = -2918.166499999999
$ws.Range("N132").Value = -11550.7367

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
$ws.Range("H132").Value = 3272.3225
$ws.Range("I132").Value = 2542.5625
$ws.Range("J132").Value = 4050.7334
$ws.Range("K132").Value = 7627.6875
$ws.Range("L132").Value = 12152.2002
$ws.Range("M132").Value = -5097.6875
$ws.Range("N132").Value = -17212.2002
$ws.Range("H136").Value = 17078.953
$ws.Range("I136").Value = 2655.76
$ws.Range("J136").Value = 37111.168
$ws.Range("K136").Value = 7967.280000000001
$ws.Range("L136").Value = 111333.504
$ws.Range("M136").Value = -5417.280000000001
$ws.Range("N136").Value = -116433.504

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10533.333
$ws.Range("I62").Value = 7755.5557
$ws.Range("J62").Value = 12200
$ws.Range("K62").Value = 7755.5557
$ws.Range("L62").Value = 12200
$ws.Range("M62").Value = -7131.5557
$ws.Range("N62").Value = -13448
$ws.Range("H65").Value = 10533.333
$ws.Range("I65").Value = 7755.5557
$ws.Range("J65").Value = 12200
$ws.Range("K65").Value = 38777.7785
$ws.Range("L65").Value = 61000
$ws.Range("M65").Value = -35657.7785
$ws.Range("N65").Value = -67240
$ws.Range("H99").Value = 39624.75
$ws.Range("I99").Value = 36166.668
$ws.Range("K99").Value = 36166.668
$ws.Range("M99").Value = -33171.668
$ws.Range("H107").Value = 745
$ws.Range("I107").Value = 744.2857
$ws.Range("K107").Value = 2232.8571
$ws.Range("M107").Value = -312.8571000000002
$ws.Range("H109").Value = 144333.33
$ws.Range("J109").Value = 144333.33
$ws.Range("L109").Value = 144333.33
$ws.Range("N109").Value = -147107.33
$ws.Range("H122").Value = 3686.1936
$ws.Range("I122").Value = 1707.5652
$ws.Range("J122").Value = 9374.75
$ws.Range("K122").Value = 5122.6956
$ws.Range("L122").Value = 28124.25
$ws.Range("M122").Value = -2672.6956
$ws.Range("N122").Value = -33024.25
$ws.Range("H124").Value = 73333.336
$ws.Range("J124").Value = 73333.336
$ws.Range("L124").Value = 73333.336
$ws.Range("N124").Value = -83153.336
$ws.Range("H132").Value = 3004.3674
$ws.Range("I132").Value = 2247.465
$ws.Range("K132").Value = 6742.395
$ws.Range("M132").Value = -4212.395

Write-Output "Applied scheduled market-data refresh across all sheets."